$wb = $excel.ActiveWorkbook

# --- Sheet1 (Statistics): update/extend rows 2-42 ---
$ws1 = $wb.Worksheets.Item("Statistics")

$ws1.Range("A2").Value = "2024-08-30 17:46:51"
$ws1.Range("B2").Value = 37.77004115973082
$ws1.Range("C2").Value = 5

$ws1.Range("A3").Value = "2024-08-30 17:46:53"
$ws1.Range("B3").Value = 40.73859557846516
$ws1.Range("C3").Value = 8

$ws1.Range("A4").Value = "2024-08-30 17:46:55"
$ws1.Range("B4").Value = 42.49414241087165
$ws1.Range("C4").Value = 12

$ws1.Range("A5").Value = "2024-08-30 17:46:57"
$ws1.Range("B5").Value = 41.58947428792983
$ws1.Range("C5").Value = 14

$ws1.Range("A6").Value = "2024-08-30 17:46:59"
$ws1.Range("B6").Value = 40.83009347718625
$ws1.Range("C6").Value = 17

$ws1.Range("A7").Value = "2024-08-30 17:47:01"
$ws1.Range("B7").Value = 38.75166211219071
$ws1.Range("C7").Value = 21

$ws1.Range("A8").Value = "2024-08-30 17:47:03"
$ws1.Range("B8").Value = 32.35722412824396
$ws1.Range("C8").Value = 23

$ws1.Range("A9").Value = "2024-08-30 17:47:05"
$ws1.Range("B9").Value = 26.6200271275341
$ws1.Range("C9").Value = 26

$ws1.Range("A10").Value = "2024-08-30 17:47:07"
$ws1.Range("B10").Value = 26.19534311713419
$ws1.Range("C10").Value = 29

$ws1.Range("A11").Value = "2024-08-30 17:47:09"
$ws1.Range("B11").Value = 24.63855276244706
$ws1.Range("C11").Value = 31

$ws1.Range("A12").Value = "2024-08-30 17:47:11"
$ws1.Range("B12").Value = 19.40124352172453
$ws1.Range("C12").Value = 32

$ws1.Range("A13").Value = "2024-08-30 17:47:13"
$ws1.Range("B13").Value = 20.73606382158988
$ws1.Range("C13").Value = 34

$ws1.Range("A14").Value = "2024-08-30 17:47:15"
$ws1.Range("B14").Value = 22.56519589868218
$ws1.Range("C14").Value = 39

$ws1.Range("A15").Value = "2024-08-30 17:47:17"
$ws1.Range("B15").Value = 19.38998777883178
$ws1.Range("C15").Value = 40

$ws1.Range("A16").Value = "2024-08-30 17:47:19"
$ws1.Range("B16").Value = 16.63371161896305
$ws1.Range("C16").Value = 41

$ws1.Range("A17").Value = "2024-08-30 17:47:21"
$ws1.Range("B17").Value = 16.45943976444216
$ws1.Range("C17").Value = 40

$ws1.Range("A18").Value = "2024-08-30 17:47:23"
$ws1.Range("B18").Value = 12.27585630456286
$ws1.Range("C18").Value = 39

$ws1.Range("A19").Value = "2024-08-30 17:47:25"
$ws1.Range("B19").Value = 14.20376145335394
$ws1.Range("C19").Value = 41

$ws1.Range("A20").Value = "2024-08-30 17:47:27"
$ws1.Range("B20").Value = 13.94596721012559
$ws1.Range("C20").Value = 41

$ws1.Range("A21").Value = "2024-08-30 17:47:29"
$ws1.Range("B21").Value = 13.53278244398157
$ws1.Range("C21").Value = 40

$ws1.Range("A22").Value = "2024-08-30 17:47:31"
$ws1.Range("B22").Value = 10.6247001290159
$ws1.Range("C22").Value = 38

$ws1.Range("A23").Value = "2024-08-30 17:47:33"
$ws1.Range("B23").Value = 12.77179836720124
$ws1.Range("C23").Value = 39

$ws1.Range("A24").Value = "2024-08-30 17:47:35"
$ws1.Range("B24").Value = 8.126224276215465
$ws1.Range("C24").Value = 35

$ws1.Range("A25").Value = "2024-08-30 17:47:37"
$ws1.Range("B25").Value = 11.7039561148478
$ws1.Range("C25").Value = 37

$ws1.Range("A26").Value = "2024-08-30 17:47:39"
$ws1.Range("B26").Value = 12.69294224291906
$ws1.Range("C26").Value = 37

$ws1.Range("A27").Value = "2024-08-30 17:47:41"
$ws1.Range("B27").Value = 12.21244462341455
$ws1.Range("C27").Value = 38

$ws1.Range("A28").Value = "2024-08-30 17:47:43"
$ws1.Range("B28").Value = 12.40883255124898
$ws1.Range("C28").Value = 38

$ws1.Range("A29").Value = "2024-08-30 17:47:45"
$ws1.Range("B29").Value = 12.35205869896399
$ws1.Range("C29").Value = 39

$ws1.Range("A30").Value = "2024-08-30 17:47:47"
$ws1.Range("B30").Value = 16.30887470148119
$ws1.Range("C30").Value = 41

$ws1.Range("A31").Value = "2024-08-30 17:47:49"
$ws1.Range("B31").Value = 15.07645122862061
$ws1.Range("C31").Value = 40

$ws1.Range("A32").Value = "2024-08-30 17:47:51"
$ws1.Range("B32").Value = 14.70820794473461
$ws1.Range("C32").Value = 39

$ws1.Range("A33").Value = "2024-08-30 17:47:53"
$ws1.Range("B33").Value = 15.79663702859301
$ws1.Range("C33").Value = 41

$ws1.Range("A34").Value = "2024-08-30 17:47:55"
$ws1.Range("B34").Value = 15.47986472960563
$ws1.Range("C34").Value = 41

$ws1.Range("A35").Value = "2024-08-30 17:47:57"
$ws1.Range("B35").Value = 14.21622389460384
$ws1.Range("C35").Value = 41

$ws1.Range("A36").Value = "2024-08-30 17:47:59"
$ws1.Range("B36").Value = 11.26183079164776
$ws1.Range("C36").Value = 38

$ws1.Range("A37").Value = "2024-08-30 17:48:01"
$ws1.Range("B37").Value = 13.19390537128629
$ws1.Range("C37").Value = 39

$ws1.Range("A38").Value = "2024-08-30 17:48:03"
$ws1.Range("B38").Value = 9.157998235261848
$ws1.Range("C38").Value = 41

$ws1.Range("A39").Value = "2024-08-30 17:48:05"
$ws1.Range("B39").Value = 7.42380171494699
$ws1.Range("C39").Value = 40

$ws1.Range("A40").Value = "2024-08-30 17:48:07"
$ws1.Range("B40").Value = 7.227210829806319
$ws1.Range("C40").Value = 42

$ws1.Range("A41").Value = "2024-08-30 17:48:09"
$ws1.Range("B41").Value = 4.207789989153946
$ws1.Range("C41").Value = 42

$ws1.Range("A42").Value = "2024-08-30 17:48:11"
$ws1.Range("B42").Value = 2.969678404867152
$ws1.Range("C42").Value = 41

# --- Sheet2 (Accidents): shrink to a single data row, with new values ---
$ws2 = $wb.Worksheets.Item("Accidents")
$ws2.Rows("3:5").Delete() | Out-Null

$ws2.Range("A2").Value = "2024-08-30 17:47:34"
$ws2.Range("B2").Value = "Car and Car"
$ws2.Range("C2").Value = "28.76 and 23.27"
$ws2.Range("D2").Value = 1
